# Update the cryptocurrency price/volume table with the latest scraped
# values from the GitHub Actions run.
#
# The "Price" column (D) holds numeric-looking text (e.g. "215.93",
# "25.993.87") that must stay plain text, exactly like the original
# cells. Writing a plain numeric-looking string to a Range.Value lets
# Excel auto-coerce it to a real number, so for every D-column write we
# force text storage by setting NumberFormat="@" first, then restore the
# cell to the default "Normal" style afterwards so no stray style index
# is introduced (the source file has no explicit number format on these
# cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# --- Row 2: Bitcoin ---
Set-TextValue "D2" "25.993.87"
$ws.Range("E2").Value = "  +0.41%  "

# --- Row 3: Ethereum ---
Set-TextValue "D3" "1.641.45"
$ws.Range("E3").Value = "  +0.40%  "

# --- Row 5: BNB ---
Set-TextValue "D5" "215.93"
$ws.Range("E5").Value = "  +0.65%  "

# --- Row 6: XRP ---
$ws.Range("E6").Value = "  -0.03%  "

# --- Row 7: USDC ---
$ws.Range("E7").Value = "  +0.39%  "

# --- Row 8: Cardano ---
$ws.Range("E8").Value = "  +0.41%  "

# --- Row 9: Dogecoin ---
$ws.Range("E9").Value = "  +0.70%  "

# --- Row 10: Solana ---
Set-TextValue "D10" "19.51"

# --- Row 11: TRON ---
$ws.Range("E11").Value = "  +0.33%  "

# --- Row 12: Polkadot ---
$ws.Range("E12").Value = "  +0.54%  "

# --- Row 13: WrappedEther ---
Set-TextValue "D13" "1.646.59"
$ws.Range("E13").Value = "  +0.88%  "

# --- Row 14: Polygon ---
Set-TextValue "D14" "0.543"
$ws.Range("E14").Value = "  +0.02%  "

# --- Row 15: ShibaInu ---
$ws.Range("E15").Value = "  +0.90%  "

# --- Row 16: Litecoin ---
Set-TextValue "D16" "63.36"

# --- Row 17: WrappedBTC ---
Set-TextValue "D17" "26.022.87"
$ws.Range("E17").Value = "  +0.48%  "

# --- Row 18: Dai ---
$ws.Range("E18").Value = "  +0.38%  "

# --- Row 19: BitcoinCash ---
Set-TextValue "D19" "194.08"
$ws.Range("E19").Value = "  +0.27%  "

# --- Row 20: Uniswap ---
$ws.Range("E20").Value = "  -0.78%  "

# --- Row 21: Avalanche ---
$ws.Range("E21").Value = "  -0.10%  "

# --- Row 22: Chainlink ---
$ws.Range("E22").Value = "  -1.06%  "

# --- Row 23: Stellar ---
$ws.Range("E23").Value = "  +4.53%  "

# --- Row 24: Toncoin ---
$ws.Range("E24").Value = "  -1.39%  "

# --- Row 25: BinanceUSD ---
Set-TextValue "D25" "1.00"
$ws.Range("E25").Value = "  +0.34%  "

# --- Row 27: Cosmos ---
$ws.Range("E27").Value = "  +0.49%  "

# --- Row 28: EthereumClassic ---
$ws.Range("E28").Value = "  +0.73%  "

# --- Row 29: PancakeSwap ---
$ws.Range("E29").Value = "  +0.47%  "

# --- Row 30: Hedera ---
$ws.Range("E30").Value = "  -0.74%  "

# --- Row 31: InternetComputer(DFINITY) ---
Set-TextValue "D31" "3.29"
$ws.Range("E31").Value = "  -0.22%  "

# --- Row 32: Filecoin ---
$ws.Range("E32").Value = "  +1.20%  "

# --- Row 33: LidoDAOToken ---
$ws.Range("E33").Value = "  -1.17%  "

# --- Row 34: HuobiToken ---
$ws.Range("E34").Value = "  +1.35%  "

# --- Row 35: ARBITRUM ---
$ws.Range("E35").Value = "  +0.26%  "

# --- Row 36: Maker ---
Set-TextValue "D36" "1.129.86"
$ws.Range("E36").Value = "  -0.85%  "

# --- Row 37: ImmutableX ---
$ws.Range("E37").Value = "  -1.16%  "

# --- Row 38: MXToken ---
$ws.Range("E38").Value = "  -0.39%  "

# --- Row 39: VeChain ---
$ws.Range("E39").Value = "  +0.20%  "

# --- Row 40: FraxShare ---
$ws.Range("E40").Value = "  +0.95%  "

# --- Row 41: Quant ---
Set-TextValue "D41" "98.96"
$ws.Range("E41").Value = "  -0.52%  "

# --- Row 42: TrustWalletToken ---
Set-TextValue "D42" "0.798"
$ws.Range("E42").Value = "  -0.16%  "

# --- Row 43: RocketPoolETH ---
Set-TextValue "D43" "1.777.94"
$ws.Range("E43").Value = "  +0.52%  "

# --- Row 44: BabyDogeCoin ---
$ws.Range("E44").Value = "  +4.71%  "

# --- Row 45: Aave ---
Set-TextValue "D45" "56.55"
$ws.Range("E45").Value = "  +0.42%  "

# --- Rows 46 & 47: RenderToken and Cronos swap positions ---
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D46" "1.49"
$ws.Range("E46").Value = "  +3.29%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D47" "0.0522"
$ws.Range("E47").Value = "  -1.33%  "

# --- Row 48: EnergySwap ---
$ws.Range("E48").Value = "  +0.81%  "

# --- Row 49: Mantle ---
$ws.Range("E49").Value = "  -0.19%  "

# --- Row 50: USDD ---
$ws.Range("E50").Value = "  +0.22%  "

# --- Row 51: Algorand ---
$ws.Range("E51").Value = "  -0.81%  "
